$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.04297822947448
$ws.Range("C2").Value = 8.908426769806731
$ws.Range("D2").Value = 6.03857093891686
$ws.Range("E2").Value = 12.59899722849443
$ws.Range("F2").Value = 28.63153133745537
$ws.Range("I2").Value = 25.91616189452966
$ws.Range("K2").Value = 8.51000187259174
$ws.Range("L2").Value = 10.1822836943188
$ws.Range("M2").Value = 13.6091953609692
$ws.Range("O2").Value = 25.73360879268861

$ws.Range("B3").Value = 10.79659435792799
$ws.Range("C3").Value = 8.874182159500284
$ws.Range("D3").Value = 5.990179931017917
$ws.Range("E3").Value = 12.62948369414344
$ws.Range("F3").Value = 28.68832341867683
$ws.Range("I3").Value = 26.01196085566699
$ws.Range("K3").Value = 8.322947988888183
$ws.Range("L3").Value = 10.19043898025132
$ws.Range("M3").Value = 13.57305300077471
$ws.Range("O3").Value = 25.81635565935152

$ws.Range("B4").Value = 10.64411892195424
$ws.Range("C4").Value = 8.852898917258008
$ws.Range("D4").Value = 5.959787623883936
$ws.Range("E4").Value = 12.64984224821202
$ws.Range("F4").Value = 28.72990042667999
$ws.Range("I4").Value = 26.07520621252327
$ws.Range("K4").Value = 8.207052698407322
$ws.Range("L4").Value = 10.196834249862
$ws.Range("M4").Value = 13.55267232263872
$ws.Range("O4").Value = 25.87207355515969

$ws.Range("B5").Value = 10.58177202474388
$ws.Range("C5").Value = 8.844163446659882
$ws.Range("D5").Value = 5.947236017612798
$ws.Range("E5").Value = 12.65855117777669
$ws.Range("F5").Value = 28.74852682541448
$ws.Range("I5").Value = 26.10209152440321
$ws.Range("K5").Value = 8.159627079494948
$ws.Range("L5").Value = 10.19978990835933
$ws.Range("M5").Value = 13.54482817205539
$ws.Range("O5").Value = 25.89601247502327

$ws.Range("B6").Value = 10.57140921805451
$ws.Range("C6").Value = 8.842709225808324
$ws.Range("D6").Value = 5.945141909668259
$ws.Range("E6").Value = 12.66002222161827
$ws.Range("F6").Value = 28.75172130437065
$ws.Range("I6").Value = 26.10662298444678
$ws.Range("K6").Value = 8.151742118363218
$ws.Range("L6").Value = 10.20030181860284
$ws.Range("M6").Value = 13.54355367049494
$ws.Range("O6").Value = 25.90006196588904

$ws.Range("B7").Value = 10.64327882869927
$ws.Range("C7").Value = 8.852781357276275
$ws.Range("D7").Value = 5.959619016684732
$ws.Range("E7").Value = 12.64995802862988
$ws.Range("F7").Value = 28.73014481691929
$ws.Range("I7").Value = 26.07556429382877
$ws.Range("K7").Value = 8.206413811330169
$ws.Range("L7").Value = 10.19687269499209
$ws.Range("M7").Value = 13.55256465916673
$ws.Range("O7").Value = 25.87239141172276

$ws.Range("B8").Value = 10.95832512890254
$ws.Range("C8").Value = 8.896673885170205
$ws.Range("D8").Value = 6.022029543047414
$ws.Range("E8").Value = 12.60916883306042
$ws.Range("F8").Value = 28.6497196790947
$ws.Range("I8").Value = 25.94827494846884
$ws.Range("K8").Value = 8.445760528502083
$ws.Range("L8").Value = 10.18480792100555
$ws.Range("M8").Value = 13.59636119950379
$ws.Range("O8").Value = 25.7611197638856

$ws.Range("B9").Value = 11.5629375680847
$ws.Range("C9").Value = 8.980612265258522
$ws.Range("D9").Value = 6.138827366896582
$ws.Range("E9").Value = 12.54217726433015
$ws.Range("F9").Value = 28.54533338721687
$ws.Range("I9").Value = 25.7337770199134
$ws.Range("K9").Value = 8.904112479171697
$ws.Range("L9").Value = 10.17213579006331
$ws.Range("M9").Value = 13.69634651108435
$ws.Range("O9").Value = 25.58194735001221

$ws.Range("B10").Value = 11.99441598483743
$ws.Range("C10").Value = 9.040853725767345
$ws.Range("D10").Value = 6.220979651241355
$ws.Range("E10").Value = 12.50086098351656
$ws.Range("F10").Value = 28.50128037746873
$ws.Range("I10").Value = 25.59760700412099
$ws.Range("K10").Value = 9.230710984420641
$ws.Range("L10").Value = 10.16948654731728
$ws.Range("M10").Value = 13.77802261617274
$ws.Range("O10").Value = 25.47418913155712

$ws.Range("B11").Value = 12.18704041500555
$ws.Range("C11").Value = 9.06792202809431
$ws.Range("D11").Value = 6.257506367075069
$ws.Range("E11").Value = 12.48377663380294
$ws.Range("F11").Value = 28.48834414320971
$ws.Range("I11").Value = 25.54031461227211
$ws.Range("K11").Value = 9.376421774492183
$ws.Range("L11").Value = 10.16971840152527
$ws.Range("M11").Value = 13.81687589485268
$ws.Range("O11").Value = 25.43036949716309

$ws.Range("B12").Value = 12.25938774354633
$ws.Range("C12").Value = 9.078121432154097
$ws.Range("D12").Value = 6.271212279741432
$ws.Range("E12").Value = 12.47755286245489
$ws.Range("F12").Value = 28.48446746422086
$ws.Range("I12").Value = 25.51928884190439
$ws.Range("K12").Value = 9.431136848611516
$ws.Range("L12").Value = 10.17001190338481
$ws.Range("M12").Value = 13.83182492675009
$ws.Range("O12").Value = 25.41452524178055

$ws.Range("B13").Value = 12.24383398330517
$ws.Range("C13").Value = 9.075927105727187
$ws.Range("D13").Value = 6.268266138781129
$ws.Range("E13").Value = 12.47888234002376
$ws.Range("F13").Value = 28.48525692147311
$ws.Range("I13").Value = 25.52378732582814
$ws.Range("K13").Value = 9.419374320355303
$ws.Range("L13").Value = 10.16993955954614
$ws.Range("M13").Value = 13.82859501478134
$ws.Range("O13").Value = 25.41790423754893

$ws.Range("B14").Value = 12.19300478520256
$ws.Range("C14").Value = 9.068762171018363
$ws.Range("D14").Value = 6.258636510566619
$ws.Range("E14").Value = 12.48325967784967
$ws.Range("F14").Value = 28.48800472352998
$ws.Range("I14").Value = 25.53857138559162
$ws.Range("K14").Value = 9.380932764264323
$ws.Range("L14").Value = 10.16973843113989
$ws.Range("M14").Value = 13.81810107549035
$ws.Range("O14").Value = 25.42905095546769

$ws.Range("B15").Value = 12.16179092427907
$ws.Range("C15").Value = 9.064366755623228
$ws.Range("D15").Value = 6.252721549824491
$ws.Range("E15").Value = 12.48597291252325
$ws.Range("F15").Value = 28.48982092971019
$ws.Range("I15").Value = 25.54771427247388
$ws.Range("K15").Value = 9.357324523621054
$ws.Range("L15").Value = 10.16964199371841
$ws.Range("M15").Value = 13.81170374187678
$ws.Range("O15").Value = 25.43597626524542

$ws.Range("B16").Value = 11.98174788037804
$ws.Range("C16").Value = 9.039077778526664
$ws.Range("D16").Value = 6.218575159949367
$ws.Range("E16").Value = 12.50201188408263
$ws.Range("F16").Value = 28.50226879049661
$ws.Range("I16").Value = 25.60144486344398
$ws.Range("K16").Value = 9.221126429140329
$ws.Range("L16").Value = 10.16950023317232
$ws.Range("M16").Value = 13.77551695435379
$ws.Range("O16").Value = 25.47715761910206

$ws.Range("B17").Value = 11.87030958624783
$ws.Range("C17").Value = 9.023476028716598
$ws.Range("D17").Value = 6.19740786870825
$ws.Range("E17").Value = 12.51228920799178
$ws.Range("F17").Value = 28.51172511029242
$ws.Range("I17").Value = 25.63559880211345
$ws.Range("K17").Value = 9.136803147433238
$ws.Range("L17").Value = 10.16978067921021
$ws.Range("M17").Value = 13.75374687163811
$ws.Range("O17").Value = 25.50375402767994

$ws.Range("B18").Value = 11.80587282400429
$ws.Range("C18").Value = 9.014470852619969
$ws.Range("D18").Value = 6.185153751258059
$ws.Range("E18").Value = 12.51836148360771
$ws.Range("F18").Value = 28.51783277658995
$ws.Range("I18").Value = 25.65568107526719
$ws.Range("K18").Value = 9.088036206367684
$ws.Range("L18").Value = 10.17007730344557
$ws.Range("M18").Value = 13.74138566125196
$ws.Range("O18").Value = 25.51954095259224

$ws.Range("B19").Value = 11.78399937725238
$ws.Range("C19").Value = 9.011416526620122
$ws.Range("D19").Value = 6.18099124436279
$ws.Range("E19").Value = 12.52044511750799
$ws.Range("F19").Value = 28.52001554283626
$ws.Range("I19").Value = 25.66255576594283
$ws.Range("K19").Value = 9.071480417109363
$ws.Range("L19").Value = 10.17020100371155
$ws.Range("M19").Value = 13.73722815248252
$ws.Range("O19").Value = 25.52497013699667

$ws.Range("B20").Value = 11.88220813232203
$ws.Range("C20").Value = 9.025140135814327
$ws.Range("D20").Value = 6.19966939800208
$ws.Range("E20").Value = 12.51117850528045
$ws.Range("F20").Value = 28.51064926164579
$ws.Range("I20").Value = 25.63191773843239
$ws.Range("K20").Value = 9.145807459920841
$ws.Range("L20").Value = 10.16973682676211
$ws.Range("M20").Value = 13.75604779458362
$ws.Range("O20").Value = 25.50087213492165

$ws.Range("B21").Value = 12.20795123940098
$ws.Range("C21").Value = 9.070868083807635
$ws.Range("D21").Value = 6.261468419629221
$ws.Range("E21").Value = 12.48196728169863
$ws.Range("F21").Value = 28.48716988988894
$ws.Range("I21").Value = 25.53421077039987
$ws.Range("K21").Value = 9.39223692074137
$ws.Range("L21").Value = 10.16979193254944
$ws.Range("M21").Value = 13.82117705699556
$ws.Range("O21").Value = 25.4257565470381

$ws.Range("B22").Value = 12.417344283403
$ws.Range("C22").Value = 9.100456837254578
$ws.Range("D22").Value = 6.301121439988389
$ws.Range("E22").Value = 12.46430805188085
$ws.Range("F22").Value = 28.47778152708342
$ws.Range("I22").Value = 25.47425696134854
$ws.Range("K22").Value = 9.550576240012875
$ws.Range("L22").Value = 10.17102645477328
$ws.Range("M22").Value = 13.86511524151966
$ws.Range("O22").Value = 25.38103224291591

$ws.Range("B23").Value = 12.30592907426319
$ws.Range("C23").Value = 9.08469276748435
$ws.Range("D23").Value = 6.280026696536843
$ws.Range("E23").Value = 12.47360218505519
$ws.Range("F23").Value = 28.48224721892075
$ws.Range("I23").Value = 25.50589806074654
$ws.Range("K23").Value = 9.466332070887006
$ws.Range("L23").Value = 10.17025823363539
$ws.Range("M23").Value = 13.84154177193549
$ws.Range("O23").Value = 25.40450228688665

$ws.Range("B24").Value = 11.8768299449183
$ws.Range("C24").Value = 9.024387903460191
$ws.Range("D24").Value = 6.198647224272556
$ws.Range("E24").Value = 12.51168014403427
$ws.Range("F24").Value = 28.51113356242452
$ws.Range("I24").Value = 25.63358055586037
$ws.Range("K24").Value = 9.141737504510663
$ws.Range("L24").Value = 10.16975623065489
$ws.Range("M24").Value = 13.75500706483616
$ws.Range("O24").Value = 25.50217349281432

$ws.Range("B25").Value = 11.40128735973736
$ws.Range("C25").Value = 8.958146278129279
$ws.Range("D25").Value = 6.107855067381946
$ws.Range("E25").Value = 12.55891101244201
$ws.Range("F25").Value = 28.56784812859134
$ws.Range("I25").Value = 25.78804359965015
$ws.Range("K25").Value = 8.781662443960291
$ws.Range("L25").Value = 10.17439123026906
$ws.Range("M25").Value = 13.66782606666564
$ws.Range("O25").Value = 25.626230619608
